$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F3, F6, F9
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5675
$ws1.Range("F6").Value = 95
$ws1.Range("F9").Value = 538

# Sheet "全部类型" (sheet4): update F3, F7, F11
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5675
$ws4.Range("F7").Value = 95
$ws4.Range("F11").Value = 538
